# Add "Planner's cost with 10 pliers" execution results (column I) plus the
# "already equipped with 10 pliers" annotation (J4), per commit:
# "Add execution results with pliers to the Excel file"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: "Planner's cost with 10 pliers" -------------------------
$ws.Range("I1").Value = "Planner's cost with 10 pliers"

$ws.Range("I2").Value  = 18
$ws.Range("I3").Value  = 19
$ws.Range("I4").Value  = 25
$ws.Range("I5").Value  = 12
$ws.Range("I6").Value  = 29
$ws.Range("I7").Value  = 118
$ws.Range("I8").Value  = 117
$ws.Range("I9").Value  = 128
$ws.Range("I10").Value = 103
$ws.Range("I11").Value = -1
$ws.Range("I12").Value = -1
$ws.Range("I13").Value = -1
$ws.Range("I14").Value = -1
$ws.Range("I15").Value = -1
$ws.Range("I16").Value = -1
$ws.Range("I17").Value = 429

# --- Annotation call-out next to the row that already satisfies the ---------
# --- "10 pliers" requirement (row 4) ----------------------------------------
$ws.Range("J4").Value = "alerady equipped with 10 pliers"
$ws.Range("J4").Style = "Neutral"

# H4 (same row) picks up the same "Neutral" highlight, keeping its original
# (larger) font size.
$ws.Range("H4").Style = "Neutral"
$ws.Range("H4").Font.Size = 12

# --- Drop the stale chart-helper defined names ------------------------------
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# --- Misc sheet/print metadata ----------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("F21").Select()
